$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 249171854.95952606
$ws.Range("C2").Value = 322023700.27230215
$ws.Range("D2").Value = 394875545.58508
$ws.Range("E2").Value = 467727390.8978569
$ws.Range("F2").Value = 540579236.2106316

$ws.Range("B3").Value = 644165076.8081288
$ws.Range("C3").Value = 717016922.1209049
$ws.Range("D3").Value = 789868767.4336827
$ws.Range("E3").Value = 862720612.7464597
$ws.Range("F3").Value = 935572458.0592343

$ws.Range("B4").Value = 1434504683.9378834
$ws.Range("C4").Value = 1507356529.2506595
$ws.Range("D4").Value = 1580208374.5634375
$ws.Range("E4").Value = 1653060219.876214
$ws.Range("F4").Value = 1725912065.1889887

$ws.Range("B5").Value = 2383534100.1871724
$ws.Range("C5").Value = 2456385945.4999485
$ws.Range("D5").Value = 2529237790.8127265
$ws.Range("E5").Value = 2602089636.1255035
$ws.Range("F5").Value = 2674941481.438278
